# updated the incorrect scores and added the complete 22nd may marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marksheet")

# Row 45 (Chemistry, 22-May) previously held data that actually belonged to
# Chemistry on 23-May (row 48). Replace it with the correct 22-May values.
$ws.Range("D45").Value = 61
$ws.Range("E45").Value = 51
$ws.Range("F45").Value = 6
$ws.Range("G45").Value = 4

# Row 46 (Biology, 22-May) was never filled in - add the complete 22nd May
# marks (the values that had erroneously been placed in row 45).
$ws.Range("D46").Value = 50
$ws.Range("E46").Value = 47
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 3

# Row 48 (Chemistry, 23-May) no longer has data - it was a duplicate of the
# 22-May Chemistry entry, now correctly placed in row 45. Clear it back out.
$ws.Range("D48:G48").ClearContents()

# Update the active sheet view/selection to match the edited state.
$excel.ActiveWindow.ScrollRow = 36
$ws.Range("G47").Select()
